# Apply the edits described by the commit diff to the workbook.
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
$meta.Range("B5").Value  = 'SimpleQuantity with UCUM or EDQM codes or code not used'
$meta.Range("B8").Value  = '2025-08-13T14:10:49+00:00'
$meta.Range("B12").Value = 'simple quantity datatype requiring a UCUM or EDQM code or no code (only unti)'

# --- Elements sheet: row 2 (Quantity) ---
$elements.Range("L2").Value = 'A fixed quantity (no comparator) with UCUM or EDQM code or no code'
$elements.Range("M2").Value = 'The comparator is not used on a SimpleQuantity. The code SHALL be a UCUM or EDQM code if used.'
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`nqty-3:If a code for the unit is present, the system SHALL also be present {code.empty() or system.exists()}sqty-1:The comparator is not used on a SimpleQuantity {comparator.empty()}fr-med-smpl-quant-1:system SHALL be UCUM or EDQM if code is used {code.exists() and (system = ``http://standardterms.edqm.eu`` or system = ``http://unitsofmeasure.org``)}fr-med-smpl-quant-2:system SHALL not be used if code is not used {code.empty() and system.empty()}"
$elements.Rows.Item(2).AutoFit()

# --- Elements sheet: row 8 (Quantity.system) ---
$elements.Range("F8").Value = 0
$elements.Range("H8").Value = ''
$elements.Range("M8").Value = 'The identification of the system that provides the coded form of the unit.'
$elements.Range("O8").Value = 'Need to know the system that defines the coded form of the unit.'

# --- Elements sheet: row 9 (Quantity.code) ---
$elements.Range("F9").Value = 0
$elements.Range("N9").Value = 'The preferred system is UCUM, but SNOMED CT can also be used (for customary units) or ISO 4217 for currency.  The context of use may additionally require a code from a particular system.'

# --- Unhide rows 2-9 on Elements sheet (was the autofilter-hidden detail rows) ---
$elements.Range("A2:A9").EntireRow.Hidden = $false

# --- Remove AutoFilter on Elements sheet ---
if ($elements.AutoFilterMode) {
    $elements.AutoFilterMode = $false
}

# --- Remove conditional formatting on Elements sheet ---
$elements.Cells.FormatConditions.Delete()

# --- Remove the defined name (_FilterDatabase) left over from the autofilter ---
foreach ($n in @($wb.Names)) {
    if ($n.Name -like '*_FilterDatabase*') {
        $n.Delete()
    }
}
